$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'67.958.83"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "  +2.36%  "

# Row 3
$ws.Range('D3').Value = "'3.340.58"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "  +2.55%  "

# Row 4
$ws.Range('E4').Value = "  -0.05%  "

# Row 5
$ws.Range('D5').Value = "'582.55"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "  +1.51%  "

# Row 6
$ws.Range('D6').Value = "'177.82"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "  +3.76%  "

# Row 7
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "  -0.06%  "

# Row 8
$ws.Range('D8').Value = "'0.590"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "  +2.63%  "

# Row 9
$ws.Range('D9').Value = "'3.338.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "  +2.64%  "

# Row 10
$ws.Range('E10').Value = "  +9.02%  "

# Row 11
$ws.Range('D11').Value = "'0.583"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "  +3.08%  "

# Row 12
$ws.Range('D12').Value = "'47.27"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "  +5.91%  "

# Row 13
$ws.Range('E13').Value = "  +3.91%  "

# Row 14
$ws.Range('D14').Value = "'707.59"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "  +3.30%  "

# Row 15
$ws.Range('D15').Value = "'3.889.81"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "  +2.49%  "

# Row 16
$ws.Range('D16').Value = "'8.45"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "  +3.23%  "

# Row 17
$ws.Range('D17').Value = "'68.025.85"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "  +2.08%  "

# Row 18
$ws.Range('E18').Value = "  +0.01%  "

# Row 19
$ws.Range('D19').Value = "'3.346.58"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "  +2.27%  "

# Row 20
$ws.Range('D20').Value = "'17.57"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "  +2.87%  "

# Row 21
$ws.Range('D21').Value = "'11.11"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "  +4.74%  "

# Row 22
$ws.Range('D22').Value = "'0.898"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "  +2.47%  "

# Row 23
$ws.Range('D23').Value = "'5.40"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "  +4.54%  "

# Row 24
$ws.Range('D24').Value = "'17.12"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "  +2.44%  "

# Row 25
$ws.Range('D25').Value = "'100.66"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "  +1.48%  "

# Row 26
$ws.Range('D26').Value = "'3.93"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "  +2.97%  "

# Row 27
$ws.Range('D27').Value = "'2.71"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "  +2.61%  "

# Row 28
$ws.Range('D28').Value = "'9.65"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "  +6.89%  "

# Row 29
$ws.Range('D29').Value = "'33.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "  +0.24%  "

# Row 30
$ws.Range('D30').Value = "'8.60"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "  +4.29%  "

# Row 31
$ws.Range('D31').Value = "'7.04"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "  +6.80%  "

# Row 32
$ws.Range('D32').Value = "'570.55"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "  +0.52%  "

# Row 33
$ws.Range('D33').Value = "'11.03"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "  +3.02%  "

# Row 34
$ws.Range('E34').Value = "  +4.26%  "

# Row 35
$ws.Range('D35').Value = "'57.85"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "  +5.41%  "

# Row 36
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = "'0.999"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "  -0.01%  "

# Row 37
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = "'3.699.03"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "  -3.30%  "

# Row 38
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = "'3.42"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "  +3.22%  "

# Row 39
$ws.Range('D39').Value = "'34.82"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "  +11.90%  "

# Row 40
$ws.Range('E40').Value = "  +4.76%  "

# Row 41
$ws.Range('D41').Value = "'2.66"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "  +4.18%  "

# Row 42
$ws.Range('D42').Value = "'3.18"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "  +8.02%  "

# Row 43
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = "'3.37"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "  -1.21%  "

# Row 44
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = "'0.0₃0680"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "  +3.83%  "

# Row 45
$ws.Range('D45').Value = "'0.338"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "  +5.09%  "

# Row 46
$ws.Range('E46').Value = "  +2.68%  "

# Row 47
$ws.Range('D47').Value = "'2.68"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "  +7.50%  "

# Row 48
$ws.Range('E48').Value = "  +2.40%  "

# Row 49
$ws.Range('E49').Value = "  -0.25%  "

# Row 50
$ws.Range('D50').Value = "'1.35"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "  +0.01%  "

# Row 51
$ws.Range('D51').Value = "'131.20"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "  +1.73%  "
